$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.693.20'
$ws.Range("E2").Value = '  +0.90%  '
$ws.Range("D3").Value = '2.999.03'
$ws.Range("E3").Value = '  +2.98%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '381.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.73%  '
$ws.Range("E7").Value = '  +1.50%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.600'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.66'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.04%  '
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0846'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.74'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.66%  '
$ws.Range("D14").Value = '3.478.96'
$ws.Range("E14").Value = '  +3.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.61%  '
$ws.Range("D16").Value = '3.005.71'
$ws.Range("E16").Value = '  +3.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.972'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.70%  '
$ws.Range("D18").Value = '51.755.28'
$ws.Range("E18").Value = '  +1.15%  '
$ws.Range("E19").Value = '  +2.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.60%  '
$ws.Range("D22").Value = '0.0₃0961'
$ws.Range("E22").Value = '  +1.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.38'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.171'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.23'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +17.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '26.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.107'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.94'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.90'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '51.51'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.09'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0458'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.12%  '
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.56'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.64'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.86'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.116'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '124.09'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.52%  '
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.276'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +16.34%  '
$ws.Range("D47").Value = '2.061.55'
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.28'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.86%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0356'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +13.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.21'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.83%  '
